$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '281.51'; E = '1.41%' }
    @{ Row = 3; D = '27.61'; E = '1.71%' }
    @{ Row = 4; D = '4.863'; E = '-1.03%' }
    @{ Row = 5; D = '0.06518'; E = '1.64%' }
    @{ Row = 6; D = '7.092'; E = '2.02%' }
    @{ Row = 7; D = '1.281'; E = '2.98%' }
    @{ Row = 8; D = '0.9149'; E = '3.50%' }
    @{ Row = 9; D = '0.1551'; E = '1.82%' }
    @{ Row = 10; D = '0.06419'; E = '28.05%' }
    @{ Row = 11; D = '0.07542'; E = '0.37%' }
    @{ Row = 12; D = '0.02920'; E = '0.33%' }
    @{ Row = 13; D = '0.09000'; E = '-0.05%' }
    @{ Row = 14; D = '0.001604'; E = '2.28%' }
    @{ Row = 15; D = '0.0006430'; E = '0.28%' }
    @{ Row = 16; D = '0.006016'; E = '3.04%' }
    @{ Row = 17; D = '3.497'; E = '1.07%' }
    @{ Row = 18; D = '3.315'; E = '0.02%' }
    @{ Row = 19; D = '2.230'; E = '-1.84%' }
    @{ Row = 20; D = $null; E = '0.29%' }
    @{ Row = 21; D = '0.1349'; E = '0.97%' }
    @{ Row = 22; D = '3.897'; E = '-0.22%' }
    @{ Row = 23; D = '0.04408'; E = '-0.70%' }
    @{ Row = 24; D = '0.1498'; E = '8.53%' }
    @{ Row = 25; D = '0.001172'; E = '-0.10%' }
    @{ Row = 26; D = '0.004347'; E = '12.57%' }
    @{ Row = 28; D = '0.0001176'; E = '-2.10%' }
    @{ Row = 29; D = '0.0001652'; E = '-14.68%' }
    @{ Row = 40; D = '0.04116'; E = '-0.65%' }
    @{ Row = 41; D = $null; E = '19.30%' }
    @{ Row = 42; D = '0.006598'; E = '-2.95%' }
    @{ Row = 43; D = '0.002072'; E = '-13.38%' }
    @{ Row = 44; D = '0.01165'; E = '-0.59%' }
    @{ Row = 45; D = '0.00005515'; E = '5.93%' }
    @{ Row = 46; D = '1.561'; E = '5.00%' }
    @{ Row = 47; D = '0.01843'; E = '-8.97%' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
